$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-10: full Fgf7 -> Nrp1 sender/target cluster grid (ECs, FAPs, sCs) x (ECs, FAPs, sCs)
$rows = @(
    @("ECs", "Fgf7", "Nrp1", "ECs", 2, 0.6666666666666666, 0.245245, 0.735735, 0.0130094690177091, 0.0130094690177091, 3, 1, 52.91030366666666, 158.730911, 0.4161415425564564, 0.4161415425564564, 12.97598742273167, 116.783886804585, 0.005413780504869892, 0.005413780504869892),
    @("ECs", "Fgf7", "Nrp1", "FAPs", 2, 0.6666666666666666, 0.245245, 0.735735, 0.0130094690177091, 0.0130094690177091, 3, 1, 47.26005833333333, 141.780175, 0.3717021489810786, 0.3717021489810786, 11.59029300595833, 104.312637053625, 0.004835647590985233, 0.004835647590985233),
    @("ECs", "Fgf7", "Nrp1", "sCs", 2, 0.6666666666666666, 0.245245, 0.735735, 0.0130094690177091, 0.0130094690177091, 3, 1, 26.97460733333333, 80.923822, 0.2121563084624651, 0.2121563084624651, 6.615387575463334, 59.53848817917, 0.002760040921853975, 0.002760040921853975),
    @("FAPs", "Fgf7", "Nrp1", "ECs", 3, 1, 17.39906333333333, 52.19719000000001, 0.9229650976458578, 0.9229650976458579, 3, 1, 52.91030366666666, 158.730911, 0.4161415425564564, 0.4161415425564564, 920.5897244822322, 8285.307520340091, 0.3840841194601176, 0.3840841194601177),
    @("FAPs", "Fgf7", "Nrp1", "FAPs", 3, 1, 17.39906333333333, 52.19719000000001, 0.9229650976458578, 0.9229650976458579, 3, 1, 47.26005833333333, 141.780175, 0.3717021489810786, 0.3717021489810786, 822.2807480786944, 7400.52673270825, 0.3430681102294963, 0.3430681102294964),
    @("FAPs", "Fgf7", "Nrp1", "sCs", 3, 1, 17.39906333333333, 52.19719000000001, 0.9229650976458578, 0.9229650976458579, 3, 1, 26.97460733333333, 80.923822, 0.2121563084624651, 0.2121563084624651, 469.3329013844645, 4223.99611246018, 0.1958128679562439, 0.1958128679562439),
    @("sCs", "Fgf7", "Nrp1", "ECs", 3, 1, 1.206960666666667, 3.620882, 0.06402543333643303, 0.06402543333643303, 3, 1, 52.91030366666666, 158.730911, 0.4161415425564564, 0.4161415425564564, 63.86065538705577, 574.7458984835019, 0.0266436425914688, 0.0266436425914688),
    @("sCs", "Fgf7", "Nrp1", "FAPs", 3, 1, 1.206960666666667, 3.620882, 0.06402543333643303, 0.06402543333643303, 3, 1, 47.26005833333333, 141.780175, 0.3717021489810786, 0.3717021489810786, 57.04103151270554, 513.3692836143499, 0.02379839116059694, 0.02379839116059694),
    @("sCs", "Fgf7", "Nrp1", "sCs", 3, 1, 1.206960666666667, 3.620882, 0.06402543333643303, 0.06402543333643303, 3, 1, 26.97460733333333, 80.923822, 0.2121563084624651, 0.2121563084624651, 32.55729005011155, 293.015610451004, 0.01358339958436728, 0.01358339958436728)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
